$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3: Invalid and Absent marked
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count and Real marked
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count and Real marked
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Total Attendance Count and Real marked
$ws.Range("D6").Value = 1
$ws.Range("E6").Value = 1

# Row 7: Absent marked
$ws.Range("H7").Value = 1

# Row 8: Absent marked
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count and Real marked
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Total Attendance Count and Real marked
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11: Total Attendance Count and Real marked
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

# Row 12: Total Attendance Count and Real marked
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13: Total Attendance Count and Real marked
$ws.Range("D13").Value = 1
$ws.Range("E13").Value = 1

# Row 14: Absent marked
$ws.Range("H14").Value = 1

# Row 15: Total Attendance Count and Real marked
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 1

# Row 16: Absent marked
$ws.Range("H16").Value = 1

# Row 17: Total Attendance Count and Real marked
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 1

# Row 18: Absent marked
$ws.Range("H18").Value = 1
